$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (e.g. "211.00", "0.0584") are preserved verbatim as text, matching
# the workbook convention where Price/Volume columns are inline strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.499.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.555.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.27"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0584"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0893"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.776.37"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.556.55"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.474.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.512"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.45"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0673"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.91"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.77"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.79%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0466"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.55%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.394.46"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.01"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.04"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.49"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.87%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.520"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.778"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0466"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.30%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.71"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.18%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.35"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.690.87"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.869"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.50%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "43.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.22%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.47"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.43%  "
